$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row total right-answer marks (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row total marks (B12): 27 -> 45
$ws.Range("B12").Value = 45

# Update the Corr/Total marks summary text (E12): "23/84" -> "45/140"
$ws.Range("E12").Value = "45/140"
